$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Federal Corporate Taxes Contribution (current)
$ws.Range("H6").Value = 0.1067
$ws.Range("I6").Value = -0.0452
$ws.Range("J6").Value = 0.295
$ws.Range("K6").Value = 0.2573
$ws.Range("L6").Value = 0.3088
$ws.Range("M6").Value = 0.2796
$ws.Range("N6").Value = 0.1561
$ws.Range("O6").Value = -0.1221
$ws.Range("P6").Value = -0.2406
$ws.Range("Q6").Value = -0.226
$ws.Range("R6").Value = -0.2871
$ws.Range("S6").Value = -0.2221
$ws.Range("T6").Value = -0.1724
$ws.Range("U6").Value = -0.035
$ws.Range("V6").Value = -0.1145
$ws.Range("W6").Value = -0.4321

# Row 16: Fiscal Impact (current)
$ws.Range("H16").Value = -2.2874
$ws.Range("I16").Value = -0.5607
$ws.Range("J16").Value = 0.0638
$ws.Range("K16").Value = -0.1884
$ws.Range("L16").Value = 0.5287
$ws.Range("M16").Value = 0.2312
$ws.Range("N16").Value = -0.2545
$ws.Range("O16").Value = -0.7614
$ws.Range("P16").Value = -0.4572
$ws.Range("Q16").Value = -0.2208
$ws.Range("R16").Value = -0.8013
$ws.Range("S16").Value = -0.9297
$ws.Range("T16").Value = -0.6619
$ws.Range("U16").Value = -0.4547
$ws.Range("V16").Value = -0.1598
$ws.Range("W16").Value = -71.7788

# Row 34: Federal Corporate Taxes Contribution (difference)
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = -0.0018
$ws.Range("P34").Value = -0.0024
$ws.Range("Q34").Value = -0.0025
$ws.Range("R34").Value = -0.0032
$ws.Range("S34").Value = -0.0026
$ws.Range("T34").Value = -0.0017
$ws.Range("U34").Value = 0.0005
$ws.Range("V34").Value = -0.0004
$ws.Range("W34").Value = -0.0049

# Row 44: Fiscal Impact (difference)
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = -0.0876
$ws.Range("P44").Value = -0.0942
$ws.Range("Q44").Value = -0.1018
$ws.Range("R44").Value = -0.1151
$ws.Range("S44").Value = -0.1102
$ws.Range("T44").Value = -0.0977
$ws.Range("U44").Value = -0.0843
$ws.Range("V44").Value = -0.0689
$ws.Range("W44").Value = -1.1317
